$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - fill first for all new rows (matches shared-string insertion order)
$ws.Cells.Item(15, 1).Value = "Kilo"
$ws.Cells.Item(16, 1).Value = "At"
$ws.Cells.Item(17, 1).Value = "Gram"
$ws.Cells.Item(18, 1).Value = "Pound"
$ws.Cells.Item(19, 1).Value = "Lot"
$ws.Cells.Item(20, 1).Value = "Package"
$ws.Cells.Item(21, 1).Value = "Container"
$ws.Cells.Item(22, 1).Value = "Shooting"
$ws.Cells.Item(23, 1).Value = "Group"
$ws.Cells.Item(24, 1).Value = "Ton"
$ws.Cells.Item(25, 1).Value = "Pipe"
$ws.Cells.Item(26, 1).Value = "Barrel"
$ws.Cells.Item(27, 1).Value = "Reel"
$ws.Cells.Item(28, 1).Value = "Strip"
$ws.Cells.Item(29, 1).Value = "Roll"
$ws.Cells.Item(30, 1).Value = "Milligram"
$ws.Cells.Item(31, 1).Value = "Box"
$ws.Cells.Item(32, 1).Value = "Ration"

# Column B - fill next for all new rows
$ws.Cells.Item(15, 2).Value = "Kilo"
$ws.Cells.Item(16, 2).Value = "At"
$ws.Cells.Item(17, 2).Value = "Gram"
$ws.Cells.Item(18, 2).Value = "Pound"
$ws.Cells.Item(19, 2).Value = "Lot"
$ws.Cells.Item(20, 2).Value = "Package"
$ws.Cells.Item(21, 2).Value = "Container"
$ws.Cells.Item(22, 2).Value = "Shooting"
$ws.Cells.Item(23, 2).Value = "Group"
$ws.Cells.Item(24, 2).Value = "Ton"
$ws.Cells.Item(25, 2).Value = "Pipe"
$ws.Cells.Item(26, 2).Value = "Barrel"
$ws.Cells.Item(27, 2).Value = "Reel"
$ws.Cells.Item(28, 2).Value = "Strip"
$ws.Cells.Item(29, 2).Value = "Roll"
$ws.Cells.Item(30, 2).Value = "Milligram"
$ws.Cells.Item(31, 2).Value = "Box"
$ws.Cells.Item(32, 2).Value = "Ration"

# Column C - fill last for all new rows
$ws.Cells.Item(15, 3).Value = "Kilo"
$ws.Cells.Item(16, 3).Value = "Arroba"
$ws.Cells.Item(17, 3).Value = "Gramo"
$ws.Cells.Item(18, 3).Value = "Libra"
$ws.Cells.Item(19, 3).Value = "Lote"
$ws.Cells.Item(20, 3).Value = "Paquete"
$ws.Cells.Item(21, 3).Value = "Container"
$ws.Cells.Item(22, 3).Value = "Disparo"
$ws.Cells.Item(23, 3).Value = "Grupo"
$ws.Cells.Item(24, 3).Value = "Tonelada"
$ws.Cells.Item(25, 3).Value = "Pipa"
$ws.Cells.Item(26, 3).Value = "Barril"
$ws.Cells.Item(27, 3).Value = "Carrete"
$ws.Cells.Item(28, 3).Value = "Tira"
$ws.Cells.Item(29, 3).Value = "Rollo"
$ws.Cells.Item(30, 3).Value = "Miligramo"
$ws.Cells.Item(31, 3).Value = "Caja"
$ws.Cells.Item(32, 3).Value = "Ración"

# Formulas: per-cell assignment so the writer reuses the existing shared-formula group (si=0 / si=1)
$ws.Cells.Item(15, 4).Formula = '=CONCATENATE("''",A15,"'':","''",B15,"'',")'
$ws.Cells.Item(15, 5).Formula = '=CONCATENATE("''",A15,"'':","''",C15,"'',")'
$ws.Cells.Item(16, 4).Formula = '=CONCATENATE("''",A16,"'':","''",B16,"'',")'
$ws.Cells.Item(16, 5).Formula = '=CONCATENATE("''",A16,"'':","''",C16,"'',")'
$ws.Cells.Item(17, 4).Formula = '=CONCATENATE("''",A17,"'':","''",B17,"'',")'
$ws.Cells.Item(17, 5).Formula = '=CONCATENATE("''",A17,"'':","''",C17,"'',")'
$ws.Cells.Item(18, 4).Formula = '=CONCATENATE("''",A18,"'':","''",B18,"'',")'
$ws.Cells.Item(18, 5).Formula = '=CONCATENATE("''",A18,"'':","''",C18,"'',")'
$ws.Cells.Item(19, 4).Formula = '=CONCATENATE("''",A19,"'':","''",B19,"'',")'
$ws.Cells.Item(19, 5).Formula = '=CONCATENATE("''",A19,"'':","''",C19,"'',")'
$ws.Cells.Item(20, 4).Formula = '=CONCATENATE("''",A20,"'':","''",B20,"'',")'
$ws.Cells.Item(20, 5).Formula = '=CONCATENATE("''",A20,"'':","''",C20,"'',")'
$ws.Cells.Item(21, 4).Formula = '=CONCATENATE("''",A21,"'':","''",B21,"'',")'
$ws.Cells.Item(21, 5).Formula = '=CONCATENATE("''",A21,"'':","''",C21,"'',")'
$ws.Cells.Item(22, 4).Formula = '=CONCATENATE("''",A22,"'':","''",B22,"'',")'
$ws.Cells.Item(22, 5).Formula = '=CONCATENATE("''",A22,"'':","''",C22,"'',")'
$ws.Cells.Item(23, 4).Formula = '=CONCATENATE("''",A23,"'':","''",B23,"'',")'
$ws.Cells.Item(23, 5).Formula = '=CONCATENATE("''",A23,"'':","''",C23,"'',")'
$ws.Cells.Item(24, 4).Formula = '=CONCATENATE("''",A24,"'':","''",B24,"'',")'
$ws.Cells.Item(24, 5).Formula = '=CONCATENATE("''",A24,"'':","''",C24,"'',")'
$ws.Cells.Item(25, 4).Formula = '=CONCATENATE("''",A25,"'':","''",B25,"'',")'
$ws.Cells.Item(25, 5).Formula = '=CONCATENATE("''",A25,"'':","''",C25,"'',")'
$ws.Cells.Item(26, 4).Formula = '=CONCATENATE("''",A26,"'':","''",B26,"'',")'
$ws.Cells.Item(26, 5).Formula = '=CONCATENATE("''",A26,"'':","''",C26,"'',")'
$ws.Cells.Item(27, 4).Formula = '=CONCATENATE("''",A27,"'':","''",B27,"'',")'
$ws.Cells.Item(27, 5).Formula = '=CONCATENATE("''",A27,"'':","''",C27,"'',")'
$ws.Cells.Item(28, 4).Formula = '=CONCATENATE("''",A28,"'':","''",B28,"'',")'
$ws.Cells.Item(28, 5).Formula = '=CONCATENATE("''",A28,"'':","''",C28,"'',")'
$ws.Cells.Item(29, 4).Formula = '=CONCATENATE("''",A29,"'':","''",B29,"'',")'
$ws.Cells.Item(29, 5).Formula = '=CONCATENATE("''",A29,"'':","''",C29,"'',")'
$ws.Cells.Item(30, 4).Formula = '=CONCATENATE("''",A30,"'':","''",B30,"'',")'
$ws.Cells.Item(30, 5).Formula = '=CONCATENATE("''",A30,"'':","''",C30,"'',")'
$ws.Cells.Item(31, 4).Formula = '=CONCATENATE("''",A31,"'':","''",B31,"'',")'
$ws.Cells.Item(31, 5).Formula = '=CONCATENATE("''",A31,"'':","''",C31,"'',")'
$ws.Cells.Item(32, 4).Formula = '=CONCATENATE("''",A32,"'':","''",B32,"'',")'
$ws.Cells.Item(32, 5).Formula = '=CONCATENATE("''",A32,"'':","''",C32,"'',")'

# Update view state to match the final selection / scroll position
[void]$ws.Range("E15:E32").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
